$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 216, pushing the existing rows 216:229 down to 218:231.
$ws.Rows("216:217").Insert()

# Row 216 - newest weekly price entry.
$ws.Range("A216").Value = 10
$ws.Range("B216").Value = "Vega Modelo de Temuco"
$ws.Range("C216").Value = "La Araucanía"
$ws.Range("D216").Value = 44706
$ws.Range("E216").Value = 9
$ws.Range("F216").Value = 100112043
$ws.Range("G216").Value = "Pepino dulce"
$ws.Range("H216").Value = "Cultivar IV Región"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 80
$ws.Range("K216").Value = 18000
$ws.Range("L216").Value = 18000
$ws.Range("M216").Value = 18000
$ws.Range("N216").Value = "$/bandeja 18 kilos"
$ws.Range("O216").Value = "Provincia de Limarí"
$ws.Range("P216").Value = 1000
$ws.Range("Q216").Value = 18
$ws.Range("R216").Value = "Hortaliza"

# Row 217 - second newest weekly price entry.
$ws.Range("A217").Value = 10
$ws.Range("B217").Value = "Vega Modelo de Temuco"
$ws.Range("C217").Value = "La Araucanía"
$ws.Range("D217").Value = 44706
$ws.Range("E217").Value = 9
$ws.Range("F217").Value = 100112043
$ws.Range("G217").Value = "Pepino dulce"
$ws.Range("H217").Value = "Cultivar IV Región"
$ws.Range("I217").Value = "Segunda"
$ws.Range("J217").Value = 160
$ws.Range("K217").Value = 15000
$ws.Range("L217").Value = 15000
$ws.Range("M217").Value = 15000
$ws.Range("N217").Value = "$/bandeja 18 kilos"
$ws.Range("O217").Value = "Provincia de Limarí"
$ws.Range("P217").Value = 833
$ws.Range("Q217").Value = 18
$ws.Range("R217").Value = "Hortaliza"
